$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @(4, 'sd', 'Statement-non-opinion'),
    @(19, '%', 'Uninterpretable'),
    @(22, 'sd', 'Statement-non-opinion'),
    @(27, 'aa', 'Agree/Accept'),
    @(35, 'sd', 'Statement-non-opinion'),
    @(37, 'sd', 'Statement-non-opinion'),
    @(48, 'ba', 'Appreciation'),
    @(57, 'sv', 'Statement-opinion'),
    @(63, 'b', 'Acknowledge (Backchannel)'),
    @(72, 'aa', 'Agree/Accept'),
    @(82, 'aa', 'Agree/Accept'),
    @(83, 'sd', 'Statement-non-opinion'),
    @(84, 'ba', 'Appreciation'),
    @(89, 'aa', 'Agree/Accept'),
    @(94, 'aa', 'Agree/Accept'),
    @(95, 'sd', 'Statement-non-opinion'),
    @(101, 'sd', 'Statement-non-opinion'),
    @(109, 'sd', 'Statement-non-opinion'),
    @(114, 'aa', 'Agree/Accept'),
    @(119, 'b', 'Acknowledge (Backchannel)'),
    @(120, 'ba', 'Appreciation'),
    @(158, 'ba', 'Appreciation'),
    @(159, 'aa', 'Agree/Accept'),
    @(163, 'sd', 'Statement-non-opinion'),
    @(177, 'sd', 'Statement-non-opinion'),
    @(187, 'sd', 'Statement-non-opinion'),
    @(189, 'sd', 'Statement-non-opinion'),
    @(206, 'sd', 'Statement-non-opinion'),
    @(211, 'sd', 'Statement-non-opinion'),
    @(213, 'aa', 'Agree/Accept'),
    @(217, 'sv', 'Statement-opinion'),
    @(226, 'sd', 'Statement-non-opinion'),
    @(240, 'aa', 'Agree/Accept'),
    @(262, 'sv', 'Statement-opinion'),
    @(263, 'aa', 'Agree/Accept'),
    @(299, 'sd', 'Statement-non-opinion'),
    @(303, 'qy', 'Yes-No-Question'),
    @(320, 'qy', 'Yes-No-Question'),
    @(348, 'sd', 'Statement-non-opinion'),
    @(351, 'sd', 'Statement-non-opinion'),
    @(354, 'sd', 'Statement-non-opinion'),
    @(373, '%', 'Uninterpretable'),
    @(379, '%', 'Uninterpretable'),
    @(384, 'b', 'Acknowledge (Backchannel)'),
    @(388, 'sd', 'Statement-non-opinion'),
    @(398, 'ba', 'Appreciation'),
    @(399, 'b', 'Acknowledge (Backchannel)'),
    @(411, 'sd', 'Statement-non-opinion'),
    @(424, 'aa', 'Agree/Accept'),
    @(451, 'aa', 'Agree/Accept'),
    @(455, 'b', 'Acknowledge (Backchannel)'),
    @(478, 'aa', 'Agree/Accept'),
    @(480, 'sv', 'Statement-opinion'),
    @(484, 'sv', 'Statement-opinion'),
    @(488, 'aa', 'Agree/Accept'),
    @(489, 'sd', 'Statement-non-opinion'),
    @(506, 'b', 'Acknowledge (Backchannel)'),
    @(508, 'sd', 'Statement-non-opinion'),
    @(511, 'ba', 'Appreciation'),
    @(526, 'ba', 'Appreciation'),
    @(527, 'b', 'Acknowledge (Backchannel)'),
    @(543, 'sv', 'Statement-opinion'),
    @(544, 'sv', 'Statement-opinion'),
    @(545, 'fc', 'Conventional-closing'),
    @(547, 'sv', 'Statement-opinion'),
    @(548, '%', 'Uninterpretable'),
    @(549, 'aa', 'Agree/Accept'),
    @(560, 'aa', 'Agree/Accept')
)

foreach ($u in $updates) {
    $row = $u[0]
    $ws.Cells.Item($row, 9).Value = $u[1]
    $ws.Cells.Item($row, 10).Value = $u[2]
}

